# Player Performance workbook update
#  1. Insert a new "Player Info" sheet in front of the existing sheets.
#  2. Rename MATCH_CARD_LINK -> MATCH_CODE on both existing sheets and
#     replace the full scorecard URL with the bare numeric match code.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. New "Player Info" sheet (becomes the first sheet in the workbook)
# ---------------------------------------------------------------------
$info = $wb.Worksheets.Add()
$info.Name = "Player Info"

$info.Range("A1").Value = "ID"
$info.Range("B1").Value = "NAME"
$info.Range("C1").Value = "BATTING_HAND"
$info.Range("D1").Value = "BOWL_STYLE"

# Match the bold/centered/bordered header styling already used on the
# other two sheets by copying the format from one of their headers.
$battingHeader = $wb.Worksheets.Item("ODI Batting").Range("A1")
$battingHeader.Copy()
$info.Range("A1:D1").PasteSpecial(-4122)  # xlPasteFormats

# Keep the player ID as text (matches how the other sheets store
# numeric-looking codes as text) by formatting the cell before the
# numeric-looking string is assigned, so Excel doesn't coerce it.
$info.Range("A2").NumberFormat = "@"
$info.Range("A2").Value = "5945"
$info.Range("B2").Value = "Zahid Mahmood"
$info.Range("C2").Value = "Right Handed"
$info.Range("D2").Value = "Right Arm Leg Break"

# ---------------------------------------------------------------------
# 2. "ODI Batting": MATCH_CARD_LINK -> MATCH_CODE, URL -> bare code
# ---------------------------------------------------------------------
$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Range("D1").Value = "MATCH_CODE"

$battingCodes = @("4564", "4565", "4567", "4641")
$battingRange = $batting.Range("D2:D5")
$battingRange.NumberFormat = "@"
for ($i = 0; $i -lt $battingCodes.Length; $i++) {
    $batting.Cells.Item(2 + $i, 4).Value = $battingCodes[$i]
}

# ---------------------------------------------------------------------
# 3. "ODI Bowling": MATCH_CARD_LINK -> MATCH_CODE, URL -> bare code
# ---------------------------------------------------------------------
$bowling = $wb.Worksheets.Item("ODI Bowling")
$bowling.Range("B1").Value = "MATCH_CODE"

$bowlingCodes = @("4564", "4565", "4567", "4641")
$bowlingRange = $bowling.Range("B2:B5")
$bowlingRange.NumberFormat = "@"
for ($i = 0; $i -lt $bowlingCodes.Length; $i++) {
    $bowling.Cells.Item(2 + $i, 2).Value = $bowlingCodes[$i]
}

Write-Output "Player Info sheet added; match links converted to codes."
